# Append a new data row (row 13) to the "Artfynd" worksheet, matching the
# source dataset export format used for the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

$ws.Cells.Item($row, 1).Value  = 111826663
$ws.Cells.Item($row, 2).Value  = 90668
$ws.Cells.Item($row, 3).Value  = "Ovaliderad"
$ws.Cells.Item($row, 4).Value  = "NT"
$ws.Cells.Item($row, 5).Value  = 788
$ws.Cells.Item($row, 6).Value  = "Gul taggsvamp"
$ws.Cells.Item($row, 7).Value  = "Hydnellum geogenium"
$ws.Cells.Item($row, 8).Value  = "(Fr.) Banker"

$ws.Cells.Item($row, 16).Value = "Färgelanda, Dls"
$ws.Cells.Item($row, 17).Value = 334977.3169877924
$ws.Cells.Item($row, 18).Value = 6509731.962975406
$ws.Cells.Item($row, 19).Value = 25
$ws.Cells.Item($row, 20).Value = "Västra Götaland"
$ws.Cells.Item($row, 21).Value = "Färgelanda"
$ws.Cells.Item($row, 22).Value = "Dalsland"
$ws.Cells.Item($row, 23).Value = "Järbo"

# Dates are stored as plain text in this workbook's export format (not real
# Excel date serials). A leading apostrophe forces text entry so "2023-09-01"
# is kept verbatim instead of being auto-converted to a date value.
$ws.Cells.Item($row, 25).Value = "'2023-09-01"
$ws.Cells.Item($row, 26).Value = "00:00"
$ws.Cells.Item($row, 27).Value = "'2023-09-01"
$ws.Cells.Item($row, 28).Value = "00:00"

$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false
$ws.Cells.Item($row, 33).Value = $false

$ws.Cells.Item($row, 49).Value = "Christine Bryngelsson"
$ws.Cells.Item($row, 50).Value = "Christine Bryngelsson"
